$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 2876
$ws.Range("F3").Value = 1164
$ws.Range("F4").Value = 21061
$ws.Range("F6").Value = 2825
$ws.Range("F7").Value = 804
$ws.Range("F9").Value = 511
$ws.Range("F10").Value = 762
$ws.Range("F14").Value = 112
$ws.Range("F15").Value = 511
$ws.Range("F17").Value = 262
$ws.Range("F18").Value = 15
$ws.Range("F19").Value = 420
$ws.Range("F20").Value = 53
$ws.Range("F21").Value = 26
$ws.Range("F23").Value = 28
$ws.Range("F24").Value = 124

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F3").Value = 30
$ws.Range("F5").Value = 336
$ws.Range("F12").Value = 100
$ws.Range("F14").Value = 150

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6123
$ws.Range("F3").Value = 697
$ws.Range("F4").Value = 685
$ws.Range("F5").Value = 1576

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6123
$ws.Range("F3").Value = 697
$ws.Range("F4").Value = 685
$ws.Range("F5").Value = 1576
$ws.Range("F6").Value = 2876
$ws.Range("F7").Value = 1164
$ws.Range("F8").Value = 21061
$ws.Range("F10").Value = 30
$ws.Range("F13").Value = 336
$ws.Range("F14").Value = 2825
$ws.Range("F15").Value = 804
$ws.Range("F19").Value = 511
$ws.Range("F20").Value = 762
$ws.Range("F27").Value = 112
$ws.Range("F30").Value = 511
$ws.Range("F31").Value = 100
$ws.Range("F33").Value = 262
$ws.Range("F34").Value = 150
$ws.Range("F35").Value = 150
$ws.Range("F36").Value = 15
$ws.Range("F37").Value = 420
$ws.Range("F39").Value = 53
$ws.Range("F40").Value = 26
$ws.Range("F44").Value = 28
$ws.Range("F50").Value = 124
